$d = $word.ActiveDocument

# The "Field Level Authorization" table (the 6th table in the document) has
# two cells reading "RU" (FISO and SO columns on the "status" row) that need
# to become "RCU".
$table = $d.Tables.Item(6)

foreach ($row in $table.Rows) {
    foreach ($cell in $row.Cells) {
        $text = $cell.Range.Text
        $text = $text -replace "[\x07\x0d]+$", ""
        if ($text -eq "RU") {
            $cell.Range.Text = "RCU"
        }
    }
}

# That text change causes Word to recompute the (fixed) column grid for the
# table; reproduce the resulting widths (in twips -> points, since
# Column.Width is expressed in points).
$newWidthsTwips = @(2514, 754, 628, 502, 502, 754, 628, 502, 502, 628)
for ($i = 1; $i -le $newWidthsTwips.Length; $i++) {
    $table.Columns.Item($i).Width = $newWidthsTwips[$i - 1] / 20
}
